# Updates the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# symbol list sheet with freshly scraped values. Source cells are stored
# as plain text (e.g. "301.55", "-0.88%"), so each new value is written
# with a leading apostrophe to force Excel to keep it as text instead of
# re-interpreting it as a number/percentage, and the cell style is reset
# to "Normal" afterwards so no stray "Text" number-format style lingers
# on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue $ws "D2" "301.55"
Set-TextValue $ws "E2" "-0.88%"
Set-TextValue $ws "D3" "31.52"
Set-TextValue $ws "E3" "-2.90%"
Set-TextValue $ws "D4" "5.159"
Set-TextValue $ws "E4" "-2.70%"
Set-TextValue $ws "D5" "0.07410"
Set-TextValue $ws "E5" "-1.13%"
Set-TextValue $ws "D6" "2.191"
Set-TextValue $ws "E6" "45.89%"
Set-TextValue $ws "D7" "7.907"
Set-TextValue $ws "E7" "0.95%"
Set-TextValue $ws "E8" "-1.13%"
Set-TextValue $ws "D9" "0.9286"
Set-TextValue $ws "E9" "0.91%"
Set-TextValue $ws "D10" "0.1716"
Set-TextValue $ws "E10" "1.07%"
Set-TextValue $ws "D11" "0.07635"
Set-TextValue $ws "E11" "-2.89%"
Set-TextValue $ws "D12" "0.08178"
Set-TextValue $ws "E12" "1.21%"
Set-TextValue $ws "D13" "0.03030"
Set-TextValue $ws "E13" "1.15%"
Set-TextValue $ws "D14" "0.09931"
Set-TextValue $ws "E14" "0.31%"
Set-TextValue $ws "D15" "0.001500"
Set-TextValue $ws "E15" "0.67%"
Set-TextValue $ws "D16" "0.006165"
Set-TextValue $ws "E16" "0.76%"
Set-TextValue $ws "D17" "3.464"
Set-TextValue $ws "E17" "-0.11%"
Set-TextValue $ws "D18" "2.230"
Set-TextValue $ws "E18" "0.05%"
Set-TextValue $ws "D19" "0.3247"
Set-TextValue $ws "E19" "-2.43%"
Set-TextValue $ws "D20" "0.1336"
Set-TextValue $ws "E20" "0.67%"
Set-TextValue $ws "D21" "4.649"
Set-TextValue $ws "E21" "3.92%"
Set-TextValue $ws "D22" "0.04650"
Set-TextValue $ws "E22" "0.80%"
Set-TextValue $ws "D23" "0.1582"
Set-TextValue $ws "E23" "-2.35%"
Set-TextValue $ws "E24" "0.07%"
Set-TextValue $ws "D25" "0.004483"
Set-TextValue $ws "E25" "0.74%"
Set-TextValue $ws "D26" "0.0001299"
Set-TextValue $ws "E26" "-7.04%"
Set-TextValue $ws "E27" "7.61%"
Set-TextValue $ws "D39" "0.01733"
Set-TextValue $ws "E39" "-2.24%"
Set-TextValue $ws "D40" "0.04532"
Set-TextValue $ws "E40" "-0.30%"
Set-TextValue $ws "D41" "0.007132"
Set-TextValue $ws "E41" "-0.69%"
Set-TextValue $ws "D42" "0.1347"
Set-TextValue $ws "E42" "-0.08%"
Set-TextValue $ws "D43" "0.002188"
Set-TextValue $ws "E43" "-0.80%"
Set-TextValue $ws "E44" "-17.08%"
Set-TextValue $ws "D45" "0.00006283"
Set-TextValue $ws "E45" "0.94%"
Set-TextValue $ws "E46" "-46.12%"
Set-TextValue $ws "D47" "1.928"
Set-TextValue $ws "E47" "2.98%"
